# Update column G ("K" - strikeouts) values on Sheet1, rows 2-11,
# replacing the old "Strike#" derived values with the new "K" values.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$newValues = @{
    2  = 1
    3  = 0
    4  = 1
    5  = 2
    6  = 0
    7  = 1
    8  = 2
    9  = 2
    10 = 2
    11 = 2
}

foreach ($row in $newValues.Keys) {
    $ws.Range("G$row").Value = $newValues[$row]
}
